$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct existing data rows 3-24 (row 2 unchanged) ---
$ws.Cells.Item(3, 2).Value = "5002825-34.2018.8.21.0002"
$ws.Cells.Item(3, 3).Value = "9000482-94.2018.8.21.0002"
$ws.Cells.Item(3, 4).Value = "Migrado"
$ws.Cells.Item(3, 9).Value = "27/03/2018"
$ws.Cells.Item(4, 2).Value = "5010408-92.2022.8.21.0014"
$ws.Cells.Item(4, 3).Value = "5002665-02.2020.8.21.0014"
$ws.Cells.Item(4, 4).Value = "Originário"
$ws.Cells.Item(4, 9).Value = "02/12/2022"
$ws.Cells.Item(5, 2).Value = "5034684-18.2011.8.21.0001"
$ws.Cells.Item(5, 3).Value = "0111143-49.2011.8.21.0001"
$ws.Cells.Item(5, 4).Value = "Digitalizado"
$ws.Cells.Item(5, 9).Value = "15/04/2011"
$ws.Cells.Item(6, 2).Value = "5029337-62.2015.8.21.0001"
$ws.Cells.Item(6, 3).Value = "0090773-10.2015.8.21.0001"
$ws.Cells.Item(6, 4).Value = "Digitalizado"
$ws.Cells.Item(6, 9).Value = "17/04/2015"
$ws.Cells.Item(7, 2).Value = "5039775-79.2017.8.21.0001"
$ws.Cells.Item(7, 3).Value = "0039384-15.2017.8.21.0001"
$ws.Cells.Item(7, 4).Value = "Digitalizado"
$ws.Cells.Item(7, 9).Value = "14/03/2017"
$ws.Cells.Item(8, 2).Value = "5029499-57.2015.8.21.0001"
$ws.Cells.Item(8, 3).Value = "0136341-49.2015.8.21.0001"
$ws.Cells.Item(8, 4).Value = "Digitalizado"
$ws.Cells.Item(8, 9).Value = "08/06/2015"
$ws.Cells.Item(9, 2).Value = "5028728-21.2011.8.21.0001"
$ws.Cells.Item(9, 3).Value = "0344069-02.2011.8.21.0001"
$ws.Cells.Item(9, 4).Value = "Digitalizado"
$ws.Cells.Item(9, 9).Value = "19/10/2011"
$ws.Cells.Item(10, 2).Value = "5028728-21.2011.8.21.0001"
$ws.Cells.Item(10, 3).Value = "0344069-02.2011.8.21.0001"
$ws.Cells.Item(10, 4).Value = "Digitalizado"
$ws.Cells.Item(10, 9).Value = "19/10/2011"
$ws.Cells.Item(11, 2).Value = "5001806-14.2020.8.21.0134"
$ws.Cells.Item(11, 3).Value = "9000364-42.2020.8.21.0134"
$ws.Cells.Item(11, 9).Value = "06/07/2020"
$ws.Cells.Item(12, 2).Value = "5032111-07.2011.8.21.0001"
$ws.Cells.Item(12, 3).Value = "0111095-90.2011.8.21.0001"
$ws.Cells.Item(12, 4).Value = "Digitalizado"
$ws.Cells.Item(12, 5).Value = "Sem dados de processo originário 2"
$ws.Cells.Item(12, 6).Value = "Nulo"
$ws.Cells.Item(12, 9).Value = "15/04/2011"
$ws.Cells.Item(13, 2).Value = "5033921-46.2013.8.21.0001"
$ws.Cells.Item(13, 3).Value = "0037413-34.2013.8.21.0001"
$ws.Cells.Item(13, 9).Value = "13/02/2013"
$ws.Cells.Item(14, 2).Value = "5033938-82.2013.8.21.0001"
$ws.Cells.Item(14, 3).Value = "0084903-52.2013.8.21.0001"
$ws.Cells.Item(14, 9).Value = "28/03/2013"
$ws.Cells.Item(15, 2).Value = "5033938-82.2013.8.21.0001"
$ws.Cells.Item(15, 3).Value = "0084903-52.2013.8.21.0001"
$ws.Cells.Item(15, 9).Value = "28/03/2013"
$ws.Cells.Item(16, 2).Value = "5034971-78.2011.8.21.0001"
$ws.Cells.Item(16, 3).Value = "0108702-95.2011.8.21.0001"
$ws.Cells.Item(16, 9).Value = "13/04/2011"
$ws.Cells.Item(17, 2).Value = "5034971-78.2011.8.21.0001"
$ws.Cells.Item(17, 3).Value = "0108702-95.2011.8.21.0001"
$ws.Cells.Item(17, 9).Value = "13/04/2011"
$ws.Cells.Item(18, 2).Value = "5039859-80.2017.8.21.0001"
$ws.Cells.Item(18, 3).Value = "0129049-42.2017.8.21.0001"
$ws.Cells.Item(18, 9).Value = "06/08/2015"
$ws.Cells.Item(19, 2).Value = "5015478-52.2010.8.21.0001"
$ws.Cells.Item(19, 3).Value = "2626731-39.2010.8.21.0001"
$ws.Cells.Item(19, 9).Value = "06/10/2010"
$ws.Cells.Item(20, 2).Value = "5015600-65.2010.8.21.0001"
$ws.Cells.Item(20, 3).Value = "0155171-39.2010.8.21.0001"
$ws.Cells.Item(20, 9).Value = "19/01/2010"
$ws.Cells.Item(21, 2).Value = "5034891-17.2011.8.21.0001"
$ws.Cells.Item(21, 3).Value = "0243442-87.2011.8.21.0001"
$ws.Cells.Item(21, 9).Value = "28/07/2011"
$ws.Cells.Item(22, 2).Value = "5034891-17.2011.8.21.0001"
$ws.Cells.Item(22, 3).Value = "0243442-87.2011.8.21.0001"
$ws.Cells.Item(22, 9).Value = "28/07/2011"
$ws.Cells.Item(23, 2).Value = "5001001-39.2018.8.21.0067"
$ws.Cells.Item(23, 3).Value = "9000812-90.2018.8.21.0067"
$ws.Cells.Item(23, 4).Value = "Migrado"
$ws.Cells.Item(23, 5).Value = "Sem dados de processo originário 2"
$ws.Cells.Item(23, 6).Value = "Nulo"
$ws.Cells.Item(23, 9).Value = "29/11/2018"
$ws.Cells.Item(24, 2).Value = "5002047-92.2020.8.21.0067"
$ws.Cells.Item(24, 3).Value = "9000436-36.2020.8.21.0067"
$ws.Cells.Item(24, 4).Value = "Migrado"
$ws.Cells.Item(24, 5).Value = "Sem dados de processo originário 2"
$ws.Cells.Item(24, 6).Value = "Nulo"
$ws.Cells.Item(24, 7).Value = "Sem dados de processo originário 3"
$ws.Cells.Item(24, 8).Value = "Nulo"
$ws.Cells.Item(24, 9).Value = "24/09/2020"

# --- Append new rows 25-31, cloning formatting from the row above ---
$ws.Range("A24:I24").Copy() | Out-Null
$ws.Range("A25:I25").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(25, 1).Value = 23
$ws.Cells.Item(25, 2).Value = "5002099-25.2019.8.21.0067"
$ws.Cells.Item(25, 3).Value = "9000760-60.2019.8.21.0067"
$ws.Cells.Item(25, 4).Value = "Migrado"
$ws.Cells.Item(25, 5).Value = "Sem dados de processo originário 2"
$ws.Cells.Item(25, 6).Value = "Nulo"
$ws.Cells.Item(25, 7).Value = "Sem dados de processo originário 3"
$ws.Cells.Item(25, 8).Value = "Nulo"
$ws.Cells.Item(25, 9).Value = "27/06/2019"

$ws.Range("A25:I25").Copy() | Out-Null
$ws.Range("A26:I26").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(26, 1).Value = 24
$ws.Cells.Item(26, 2).Value = "5001002-24.2018.8.21.0067"
$ws.Cells.Item(26, 3).Value = "9000850-05.2018.8.21.0067"
$ws.Cells.Item(26, 4).Value = "Migrado"
$ws.Cells.Item(26, 5).Value = "Sem dados de processo originário 2"
$ws.Cells.Item(26, 6).Value = "Nulo"
$ws.Cells.Item(26, 7).Value = "Sem dados de processo originário 3"
$ws.Cells.Item(26, 8).Value = "Nulo"
$ws.Cells.Item(26, 9).Value = "04/12/2018"

$ws.Range("A26:I26").Copy() | Out-Null
$ws.Range("A27:I27").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(27, 1).Value = 25
$ws.Cells.Item(27, 2).Value = "5002100-10.2019.8.21.0067"
$ws.Cells.Item(27, 3).Value = "9000786-58.2019.8.21.0067"
$ws.Cells.Item(27, 4).Value = "Migrado"
$ws.Cells.Item(27, 5).Value = "Sem dados de processo originário 2"
$ws.Cells.Item(27, 6).Value = "Nulo"
$ws.Cells.Item(27, 7).Value = "Sem dados de processo originário 3"
$ws.Cells.Item(27, 8).Value = "Nulo"
$ws.Cells.Item(27, 9).Value = "03/07/2019"

$ws.Range("A27:I27").Copy() | Out-Null
$ws.Range("A28:I28").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(28, 1).Value = 26
$ws.Cells.Item(28, 2).Value = "5002101-92.2019.8.21.0067"
$ws.Cells.Item(28, 3).Value = "9000444-47.2019.8.21.0067"
$ws.Cells.Item(28, 4).Value = "Migrado"
$ws.Cells.Item(28, 5).Value = "Sem dados de processo originário 2"
$ws.Cells.Item(28, 6).Value = "Nulo"
$ws.Cells.Item(28, 7).Value = "Sem dados de processo originário 3"
$ws.Cells.Item(28, 8).Value = "Nulo"
$ws.Cells.Item(28, 9).Value = "31/05/2019"

$ws.Range("A28:I28").Copy() | Out-Null
$ws.Range("A29:I29").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(29, 1).Value = 27
$ws.Cells.Item(29, 2).Value = "5002102-77.2019.8.21.0067"
$ws.Cells.Item(29, 3).Value = "9000756-23.2019.8.21.0067"
$ws.Cells.Item(29, 4).Value = "Migrado"
$ws.Cells.Item(29, 5).Value = "Sem dados de processo originário 2"
$ws.Cells.Item(29, 6).Value = "Nulo"
$ws.Cells.Item(29, 7).Value = "Sem dados de processo originário 3"
$ws.Cells.Item(29, 8).Value = "Nulo"
$ws.Cells.Item(29, 9).Value = "27/06/2019"

$ws.Range("A29:I29").Copy() | Out-Null
$ws.Range("A30:I30").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(30, 1).Value = 28
$ws.Cells.Item(30, 2).Value = "5002103-62.2019.8.21.0067"
$ws.Cells.Item(30, 3).Value = "9000754-53.2019.8.21.0067"
$ws.Cells.Item(30, 4).Value = "Migrado"
$ws.Cells.Item(30, 5).Value = "Sem dados de processo originário 2"
$ws.Cells.Item(30, 6).Value = "Nulo"
$ws.Cells.Item(30, 7).Value = "Sem dados de processo originário 3"
$ws.Cells.Item(30, 8).Value = "Nulo"
$ws.Cells.Item(30, 9).Value = "27/06/2019"

$ws.Range("A30:I30").Copy() | Out-Null
$ws.Range("A31:I31").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(31, 1).Value = 29
$ws.Cells.Item(31, 2).Value = "5002104-47.2019.8.21.0067"
$ws.Cells.Item(31, 3).Value = "9000700-87.2019.8.21.0067"
$ws.Cells.Item(31, 4).Value = "Migrado"
$ws.Cells.Item(31, 5).Value = "Sem dados de processo originário 2"
$ws.Cells.Item(31, 6).Value = "Nulo"
$ws.Cells.Item(31, 7).Value = "Sem dados de processo originário 3"
$ws.Cells.Item(31, 8).Value = "Nulo"
$ws.Cells.Item(31, 9).Value = "19/06/2019"

$ws.Application.CutCopyMode = $false